# Auto-generated script: reorders the rows 114-127 in worksheet 'Artfynd'
# so that each row's full data (all populated columns) moves to the
# target row position described by the commit's row permutation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Artfynd')

# Row 114
$ws.Range('A114').Value = 111743524
$ws.Range('B114').Value = 94134
$ws.Range('C114').Value = 'Ovaliderad'
$ws.Range('D114').Value = 'NT'
$ws.Range('E114').Value = 53
$ws.Range('F114').Value = 'Vedtrappmossa'
$ws.Range('G114').Value = 'Crossocalyx hellerianus'
$ws.Range('H114').Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range('P114').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q114').Value = 338949.7235384365
$ws.Range('R114').Value = 6571040.381812023
$ws.Range('S114').Value = 5
$ws.Range('T114').Value = 'Västra Götaland'
$ws.Range('U114').Value = 'Bengtsfors'
$ws.Range('V114').Value = 'Dalsland'
$ws.Range('W114').Value = 'Vårvik'
$ws.Range('Y114').Value = "'2023-08-24"
$ws.Range('Z114').Value = '00:00'
$ws.Range('AA114').Value = "'2023-08-24"
$ws.Range('AB114').Value = '00:00'
$ws.Range('AD114').Value = $false
$ws.Range('AE114').Value = $false
$ws.Range('AG114').Value = $false
$ws.Range('AJ114').Value = 'tall'
$ws.Range('AK114').Value = 'Pinus sylvestris'
$ws.Range('AO114').Value = 'Pinus sylvestris'
$ws.Range('AW114').Value = 'Anton Larsson'
$ws.Range('AX114').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 115
$ws.Range('A115').Value = 111743554
$ws.Range('B115').Value = 88966
$ws.Range('C115').Value = 'Ovaliderad'
$ws.Range('D115').Value = 'NT'
$ws.Range('E115').Value = 5754
$ws.Range('F115').Value = 'Gultoppig fingersvamp'
$ws.Range('G115').Value = 'Ramaria testaceoflava'
$ws.Range('H115').Value = '(Bres.) Corner'
$ws.Range('P115').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q115').Value = 339577.2032005055
$ws.Range('R115').Value = 6571127.007499221
$ws.Range('S115').Value = 5
$ws.Range('T115').Value = 'Västra Götaland'
$ws.Range('U115').Value = 'Bengtsfors'
$ws.Range('V115').Value = 'Dalsland'
$ws.Range('W115').Value = 'Vårvik'
$ws.Range('Y115').Value = "'2023-08-24"
$ws.Range('Z115').Value = '00:00'
$ws.Range('AA115').Value = "'2023-08-24"
$ws.Range('AB115').Value = '00:00'
$ws.Range('AD115').Value = $false
$ws.Range('AE115').Value = $false
$ws.Range('AG115').Value = $false
$ws.Range('AW115').Value = 'Anton Larsson'
$ws.Range('AX115').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 116
$ws.Range('A116').Value = 111743517
$ws.Range('B116').Value = 73634
$ws.Range('C116').Value = 'Ovaliderad'
$ws.Range('D116').Value = 'LC'
$ws.Range('E116').Value = 6426
$ws.Range('F116').Value = 'Kattfotslav'
$ws.Range('G116').Value = 'Felipes leucopellaeus'
$ws.Range('H116').Value = '(Ach.) Frisch & G.Thor'
$ws.Range('P116').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q116').Value = 339278.3213300391
$ws.Range('R116').Value = 6571107.378548244
$ws.Range('S116').Value = 5
$ws.Range('T116').Value = 'Västra Götaland'
$ws.Range('U116').Value = 'Bengtsfors'
$ws.Range('V116').Value = 'Dalsland'
$ws.Range('W116').Value = 'Vårvik'
$ws.Range('Y116').Value = "'2023-08-24"
$ws.Range('Z116').Value = '00:00'
$ws.Range('AA116').Value = "'2023-08-24"
$ws.Range('AB116').Value = '00:00'
$ws.Range('AD116').Value = $false
$ws.Range('AE116').Value = $false
$ws.Range('AG116').Value = $false
$ws.Range('AW116').Value = 'Anton Larsson'
$ws.Range('AX116').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 117
$ws.Range('K117').Value = $null
$ws.Range('A117').Value = 111743523
$ws.Range('B117').Value = 73634
$ws.Range('C117').Value = 'Ovaliderad'
$ws.Range('D117').Value = 'LC'
$ws.Range('E117').Value = 6426
$ws.Range('F117').Value = 'Kattfotslav'
$ws.Range('G117').Value = 'Felipes leucopellaeus'
$ws.Range('H117').Value = '(Ach.) Frisch & G.Thor'
$ws.Range('P117').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q117').Value = 339009.0243061834
$ws.Range('R117').Value = 6571011.238422027
$ws.Range('S117').Value = 5
$ws.Range('T117').Value = 'Västra Götaland'
$ws.Range('U117').Value = 'Bengtsfors'
$ws.Range('V117').Value = 'Dalsland'
$ws.Range('W117').Value = 'Vårvik'
$ws.Range('Y117').Value = "'2023-08-24"
$ws.Range('Z117').Value = '00:00'
$ws.Range('AA117').Value = "'2023-08-24"
$ws.Range('AB117').Value = '00:00'
$ws.Range('AD117').Value = $false
$ws.Range('AE117').Value = $false
$ws.Range('AG117').Value = $false
$ws.Range('AW117').Value = 'Anton Larsson'
$ws.Range('AX117').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 118
$ws.Range('M118').Value = $null
$ws.Range('AJ118').Value = $null
$ws.Range('AK118').Value = $null
$ws.Range('AO118').Value = $null
$ws.Range('A118').Value = 111743549
$ws.Range('B118').Value = 96348
$ws.Range('C118').Value = 'Ovaliderad'
$ws.Range('D118').Value = 'VU'
$ws.Range('E118').Value = 220787
$ws.Range('F118').Value = 'Knärot'
$ws.Range('G118').Value = 'Goodyera repens'
$ws.Range('H118').Value = '(L.) R. Br.'
$ws.Range('K118').Value = 'blomning'
$ws.Range('P118').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q118').Value = 339495.029088294
$ws.Range('R118').Value = 6571076.196190646
$ws.Range('S118').Value = 5
$ws.Range('T118').Value = 'Västra Götaland'
$ws.Range('U118').Value = 'Bengtsfors'
$ws.Range('V118').Value = 'Dalsland'
$ws.Range('W118').Value = 'Vårvik'
$ws.Range('Y118').Value = "'2023-08-24"
$ws.Range('Z118').Value = '00:00'
$ws.Range('AA118').Value = "'2023-08-24"
$ws.Range('AB118').Value = '00:00'
$ws.Range('AD118').Value = $false
$ws.Range('AE118').Value = $false
$ws.Range('AG118').Value = $false
$ws.Range('AW118').Value = 'Anton Larsson'
$ws.Range('AX118').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 119
$ws.Range('A119').Value = 111743519
$ws.Range('B119').Value = 90666
$ws.Range('C119').Value = 'Ovaliderad'
$ws.Range('D119').Value = 'LC'
$ws.Range('E119').Value = 4364
$ws.Range('F119').Value = 'Dropptaggsvamp'
$ws.Range('G119').Value = 'Hydnellum ferrugineum'
$ws.Range('H119').Value = '(Fr.:Fr.) P. Karst.'
$ws.Range('P119').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q119').Value = 339118.4126724883
$ws.Range('R119').Value = 6571062.424656671
$ws.Range('S119').Value = 5
$ws.Range('T119').Value = 'Västra Götaland'
$ws.Range('U119').Value = 'Bengtsfors'
$ws.Range('V119').Value = 'Dalsland'
$ws.Range('W119').Value = 'Vårvik'
$ws.Range('Y119').Value = "'2023-08-24"
$ws.Range('Z119').Value = '00:00'
$ws.Range('AA119').Value = "'2023-08-24"
$ws.Range('AB119').Value = '00:00'
$ws.Range('AD119').Value = $false
$ws.Range('AE119').Value = $false
$ws.Range('AG119').Value = $false
$ws.Range('AW119').Value = 'Anton Larsson'
$ws.Range('AX119').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 120
$ws.Range('A120').Value = 111743551
$ws.Range('B120').Value = 96348
$ws.Range('C120').Value = 'Ovaliderad'
$ws.Range('D120').Value = 'VU'
$ws.Range('E120').Value = 220787
$ws.Range('F120').Value = 'Knärot'
$ws.Range('G120').Value = 'Goodyera repens'
$ws.Range('H120').Value = '(L.) R. Br.'
$ws.Range('P120').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q120').Value = 339522.8608171764
$ws.Range('R120').Value = 6571091.407599592
$ws.Range('S120').Value = 5
$ws.Range('T120').Value = 'Västra Götaland'
$ws.Range('U120').Value = 'Bengtsfors'
$ws.Range('V120').Value = 'Dalsland'
$ws.Range('W120').Value = 'Vårvik'
$ws.Range('Y120').Value = "'2023-08-24"
$ws.Range('Z120').Value = '00:00'
$ws.Range('AA120').Value = "'2023-08-24"
$ws.Range('AB120').Value = '00:00'
$ws.Range('AD120').Value = $false
$ws.Range('AE120').Value = $false
$ws.Range('AG120').Value = $false
$ws.Range('AW120').Value = 'Anton Larsson'
$ws.Range('AX120').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 121
$ws.Range('A121').Value = 111743516
$ws.Range('B121').Value = 96348
$ws.Range('C121').Value = 'Ovaliderad'
$ws.Range('D121').Value = 'VU'
$ws.Range('E121').Value = 220787
$ws.Range('F121').Value = 'Knärot'
$ws.Range('G121').Value = 'Goodyera repens'
$ws.Range('H121').Value = '(L.) R. Br.'
$ws.Range('P121').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q121').Value = 339415.5147437509
$ws.Range('R121').Value = 6571015.54325202
$ws.Range('S121').Value = 5
$ws.Range('T121').Value = 'Västra Götaland'
$ws.Range('U121').Value = 'Bengtsfors'
$ws.Range('V121').Value = 'Dalsland'
$ws.Range('W121').Value = 'Vårvik'
$ws.Range('Y121').Value = "'2023-08-24"
$ws.Range('Z121').Value = '00:00'
$ws.Range('AA121').Value = "'2023-08-24"
$ws.Range('AB121').Value = '00:00'
$ws.Range('AD121').Value = $false
$ws.Range('AE121').Value = $false
$ws.Range('AG121').Value = $false
$ws.Range('AW121').Value = 'Anton Larsson'
$ws.Range('AX121').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 122
$ws.Range('A122').Value = 111743546
$ws.Range('B122').Value = 96348
$ws.Range('C122').Value = 'Ovaliderad'
$ws.Range('D122').Value = 'VU'
$ws.Range('E122').Value = 220787
$ws.Range('F122').Value = 'Knärot'
$ws.Range('G122').Value = 'Goodyera repens'
$ws.Range('H122').Value = '(L.) R. Br.'
$ws.Range('P122').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q122').Value = 339474.5644867857
$ws.Range('R122').Value = 6571113.931964876
$ws.Range('S122').Value = 5
$ws.Range('T122').Value = 'Västra Götaland'
$ws.Range('U122').Value = 'Bengtsfors'
$ws.Range('V122').Value = 'Dalsland'
$ws.Range('W122').Value = 'Vårvik'
$ws.Range('Y122').Value = "'2023-08-24"
$ws.Range('Z122').Value = '00:00'
$ws.Range('AA122').Value = "'2023-08-24"
$ws.Range('AB122').Value = '00:00'
$ws.Range('AD122').Value = $false
$ws.Range('AE122').Value = $false
$ws.Range('AG122').Value = $false
$ws.Range('AW122').Value = 'Anton Larsson'
$ws.Range('AX122').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 123
$ws.Range('A123').Value = 111743527
$ws.Range('B123').Value = 96348
$ws.Range('C123').Value = 'Ovaliderad'
$ws.Range('D123').Value = 'VU'
$ws.Range('E123').Value = 220787
$ws.Range('F123').Value = 'Knärot'
$ws.Range('G123').Value = 'Goodyera repens'
$ws.Range('H123').Value = '(L.) R. Br.'
$ws.Range('P123').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q123').Value = 338598.1684531783
$ws.Range('R123').Value = 6571109.585305012
$ws.Range('S123').Value = 5
$ws.Range('T123').Value = 'Västra Götaland'
$ws.Range('U123').Value = 'Bengtsfors'
$ws.Range('V123').Value = 'Dalsland'
$ws.Range('W123').Value = 'Vårvik'
$ws.Range('Y123').Value = "'2023-08-24"
$ws.Range('Z123').Value = '00:00'
$ws.Range('AA123').Value = "'2023-08-24"
$ws.Range('AB123').Value = '00:00'
$ws.Range('AD123').Value = $false
$ws.Range('AE123').Value = $false
$ws.Range('AG123').Value = $false
$ws.Range('AW123').Value = 'Anton Larsson'
$ws.Range('AX123').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 124
$ws.Range('AJ124').Value = $null
$ws.Range('AK124').Value = $null
$ws.Range('AO124').Value = $null
$ws.Range('A124').Value = 111743515
$ws.Range('B124').Value = 96348
$ws.Range('C124').Value = 'Ovaliderad'
$ws.Range('D124').Value = 'VU'
$ws.Range('E124').Value = 220787
$ws.Range('F124').Value = 'Knärot'
$ws.Range('G124').Value = 'Goodyera repens'
$ws.Range('H124').Value = '(L.) R. Br.'
$ws.Range('P124').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q124').Value = 339441.7613444271
$ws.Range('R124').Value = 6571017.506567059
$ws.Range('S124').Value = 5
$ws.Range('T124').Value = 'Västra Götaland'
$ws.Range('U124').Value = 'Bengtsfors'
$ws.Range('V124').Value = 'Dalsland'
$ws.Range('W124').Value = 'Vårvik'
$ws.Range('Y124').Value = "'2023-08-24"
$ws.Range('Z124').Value = '00:00'
$ws.Range('AA124').Value = "'2023-08-24"
$ws.Range('AB124').Value = '00:00'
$ws.Range('AD124').Value = $false
$ws.Range('AE124').Value = $false
$ws.Range('AG124').Value = $false
$ws.Range('AW124').Value = 'Anton Larsson'
$ws.Range('AX124').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 125
$ws.Range('A125').Value = 111743520
$ws.Range('B125').Value = 56398
$ws.Range('C125').Value = 'Ovaliderad'
$ws.Range('D125').Value = 'NT'
$ws.Range('E125').Value = 100109
$ws.Range('F125').Value = 'Tretåig hackspett'
$ws.Range('G125').Value = 'Picoides tridactylus'
$ws.Range('H125').Value = '(Linnaeus, 1758)'
$ws.Range('M125').Value = 'färska spår'
$ws.Range('P125').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q125').Value = 339096.8530521042
$ws.Range('R125').Value = 6571013.66294401
$ws.Range('S125').Value = 5
$ws.Range('T125').Value = 'Västra Götaland'
$ws.Range('U125').Value = 'Bengtsfors'
$ws.Range('V125').Value = 'Dalsland'
$ws.Range('W125').Value = 'Vårvik'
$ws.Range('Y125').Value = "'2023-08-24"
$ws.Range('Z125').Value = '00:00'
$ws.Range('AA125').Value = "'2023-08-24"
$ws.Range('AB125').Value = '00:00'
$ws.Range('AD125').Value = $false
$ws.Range('AE125').Value = $false
$ws.Range('AG125').Value = $false
$ws.Range('AJ125').Value = 'gran'
$ws.Range('AK125').Value = 'Picea abies'
$ws.Range('AO125').Value = 'Picea abies'
$ws.Range('AW125').Value = 'Anton Larsson'
$ws.Range('AX125').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 126
$ws.Range('A126').Value = 111743526
$ws.Range('B126').Value = 90666
$ws.Range('C126').Value = 'Ovaliderad'
$ws.Range('D126').Value = 'LC'
$ws.Range('E126').Value = 4364
$ws.Range('F126').Value = 'Dropptaggsvamp'
$ws.Range('G126').Value = 'Hydnellum ferrugineum'
$ws.Range('H126').Value = '(Fr.:Fr.) P. Karst.'
$ws.Range('P126').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q126').Value = 338870.1217119552
$ws.Range('R126').Value = 6571086.774471543
$ws.Range('S126').Value = 5
$ws.Range('T126').Value = 'Västra Götaland'
$ws.Range('U126').Value = 'Bengtsfors'
$ws.Range('V126').Value = 'Dalsland'
$ws.Range('W126').Value = 'Vårvik'
$ws.Range('Y126').Value = "'2023-08-24"
$ws.Range('Z126').Value = '00:00'
$ws.Range('AA126').Value = "'2023-08-24"
$ws.Range('AB126').Value = '00:00'
$ws.Range('AD126').Value = $false
$ws.Range('AE126').Value = $false
$ws.Range('AG126').Value = $false
$ws.Range('AW126').Value = 'Anton Larsson'
$ws.Range('AX126').Value = 'Anton Larsson, Ingalill  Larsson'

# Row 127
$ws.Range('A127').Value = 111743521
$ws.Range('B127').Value = 96348
$ws.Range('C127').Value = 'Ovaliderad'
$ws.Range('D127').Value = 'VU'
$ws.Range('E127').Value = 220787
$ws.Range('F127').Value = 'Knärot'
$ws.Range('G127').Value = 'Goodyera repens'
$ws.Range('H127').Value = '(L.) R. Br.'
$ws.Range('P127').Value = 'Vårviks-Bottnane, Dls'
$ws.Range('Q127').Value = 339070.1946752003
$ws.Range('R127').Value = 6571001.989220584
$ws.Range('S127').Value = 5
$ws.Range('T127').Value = 'Västra Götaland'
$ws.Range('U127').Value = 'Bengtsfors'
$ws.Range('V127').Value = 'Dalsland'
$ws.Range('W127').Value = 'Vårvik'
$ws.Range('Y127').Value = "'2023-08-24"
$ws.Range('Z127').Value = '00:00'
$ws.Range('AA127').Value = "'2023-08-24"
$ws.Range('AB127').Value = '00:00'
$ws.Range('AD127').Value = $false
$ws.Range('AE127').Value = $false
$ws.Range('AG127').Value = $false
$ws.Range('AW127').Value = 'Anton Larsson'
$ws.Range('AX127').Value = 'Anton Larsson, Ingalill  Larsson'
